$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()
$ws.Rows("18:48").AutoFit() | Out-Null

$ws.Cells.Item(15, 3).Clear()

$ws.Cells.Item(15, 1).Value = "Criar o botao de pause"
$ws.Cells.Item(15, 4).Value = "T1.12"

$moves = @(
    @{Src="C5"; Dst="F5"},
    @{Src="D6"; Dst="F6"},
    @{Src="D7"; Dst="F7"},
    @{Src="D8"; Dst="F8"},
    @{Src="C9"; Dst="F9"},
    @{Src="C10"; Dst="F10"},
    @{Src="C11"; Dst="F11"},
    @{Src="C12"; Dst="F12"},
    @{Src="D14"; Dst="F14"}
)

foreach ($m in $moves) {
    $srcRange = $ws.Range($m.Src)
    $dstRange = $ws.Range($m.Dst)
    $srcRange.Copy($dstRange)
    $srcRange.Clear()
}

$ws.Range("D15").Select() | Out-Null

Write-Host "done"
